$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I12").Value = 'sv'
$ws.Range("J12").Value = 'Statement-opinion'
$ws.Range("I27").Value = 'b'
$ws.Range("J27").Value = 'Acknowledge (Backchannel)'
$ws.Range("I48").Value = 'aa'
$ws.Range("J48").Value = 'Agree/Accept'
$ws.Range("I54").Value = 'sd'
$ws.Range("J54").Value = 'Statement-non-opinion'
$ws.Range("I79").Value = 'ba'
$ws.Range("J79").Value = 'Appreciation'
$ws.Range("I80").Value = 'sv'
$ws.Range("J80").Value = 'Statement-opinion'
$ws.Range("I98").Value = 'sv'
$ws.Range("J98").Value = 'Statement-opinion'
$ws.Range("I99").Value = 'b'
$ws.Range("J99").Value = 'Acknowledge (Backchannel)'
$ws.Range("I102").Value = 'b'
$ws.Range("J102").Value = 'Acknowledge (Backchannel)'
$ws.Range("I108").Value = 'sv'
$ws.Range("J108").Value = 'Statement-opinion'
$ws.Range("I153").Value = 'sv'
$ws.Range("J153").Value = 'Statement-opinion'
$ws.Range("I155").Value = 'b'
$ws.Range("J155").Value = 'Acknowledge (Backchannel)'
$ws.Range("I176").Value = 'sd'
$ws.Range("J176").Value = 'Statement-non-opinion'
$ws.Range("I177").Value = 'b'
$ws.Range("J177").Value = 'Acknowledge (Backchannel)'
$ws.Range("I187").Value = 'b'
$ws.Range("J187").Value = 'Acknowledge (Backchannel)'
$ws.Range("I210").Value = 'sv'
$ws.Range("J210").Value = 'Statement-opinion'
$ws.Range("I225").Value = 'sd'
$ws.Range("J225").Value = 'Statement-non-opinion'
$ws.Range("I232").Value = 'sv'
$ws.Range("J232").Value = 'Statement-opinion'
$ws.Range("I235").Value = 'sd'
$ws.Range("J235").Value = 'Statement-non-opinion'
$ws.Range("I238").Value = 'sd'
$ws.Range("J238").Value = 'Statement-non-opinion'
$ws.Range("I253").Value = 'sv'
$ws.Range("J253").Value = 'Statement-opinion'
$ws.Range("I258").Value = 'b'
$ws.Range("J258").Value = 'Acknowledge (Backchannel)'
$ws.Range("I275").Value = 'sv'
$ws.Range("J275").Value = 'Statement-opinion'
$ws.Range("I279").Value = 'b'
$ws.Range("J279").Value = 'Acknowledge (Backchannel)'
$ws.Range("I283").Value = 'sd'
$ws.Range("J283").Value = 'Statement-non-opinion'
$ws.Range("I284").Value = 'sv'
$ws.Range("J284").Value = 'Statement-opinion'
$ws.Range("I304").Value = 'sv'
$ws.Range("J304").Value = 'Statement-opinion'
$ws.Range("I309").Value = 'sd'
$ws.Range("J309").Value = 'Statement-non-opinion'
$ws.Range("I341").Value = 'ba'
$ws.Range("J341").Value = 'Appreciation'
$ws.Range("I348").Value = 'sd'
$ws.Range("J348").Value = 'Statement-non-opinion'
$ws.Range("I352").Value = 'ba'
$ws.Range("J352").Value = 'Appreciation'
$ws.Range("I356").Value = 'b'
$ws.Range("J356").Value = 'Acknowledge (Backchannel)'
$ws.Range("I385").Value = 'sv'
$ws.Range("J385").Value = 'Statement-opinion'
$ws.Range("I386").Value = 'b'
$ws.Range("J386").Value = 'Acknowledge (Backchannel)'
